$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new header/value in column A
$ws.Range("A1").Value = "zip_code"
$ws.Range("A2").Value = 94553

# Remove the old B/C columns content and formatting entirely (no shifting,
# so column widths in <cols> and A1's own style stay untouched)
$ws.Range("B1:C2").Clear()

# Match the saved selection state
$ws.Range("A2").Select()
